{"js": "// The document has a number of paragraphs whose <w:pPr> carries two\n// <w:pStyle> entries \u2014 a leftover \"Compact\" style plus the paragraph's\n// real style (Word/Office.js only ever \"sees\" the second/last one as the\n// paragraph's effective style). This edit removes the redundant\n// \"Compact\" entry:\n//   - paragraphs that resolve to a real style (Casebook Title/Subtitle,\n//     Section Number/Title/Subtitle, Resource Number/Title/Link) keep\n//     that style \u2014 re-applying it collapses the duplicate pStyle down to\n//     a single entry;\n//   - the one paragraph whose *only* style was \"Compact\" is switched to\n//     \"Body Text\" instead.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst items = paragraphs.items;\nfor (const p of items) {\n  p.load(\"style\");\n}\nawait context.sync();\n\n// Styles that, when seen as a paragraph's *current* effective style,\n// mark it as one of the \"Compact + real style\" paragraphs targeted by\n// this edit. Re-assigning the same style collapses the stray duplicate\n// \"Compact\" pStyle entry that precedes it.\nconst keepStyle = new Set([\n  \"Casebook Title\",\n  \"Casebook Subtitle\",\n  \"Section Number\",\n  \"Section Title\",\n  \"Section Subtitle\",\n  \"Resource Number\",\n  \"Resource Title\",\n  \"Resource Link\",\n]);\n\nfor (const p of items) {\n  if (p.style === \"Compact\") {\n    // The lone paragraph whose only pStyle was \"Compact\" becomes \"Body Text\".\n    p.style = \"Body Text\";\n  } else if (keepStyle.has(p.style)) {\n    p.style = p.style;\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document has a number of paragraphs whose <w:pPr> carries two\n# <w:pStyle> entries -- a leftover \"Compact\" style plus the paragraph's\n# real style (Word only ever \"sees\" the second/last one as the\n# paragraph's effective style). This edit removes the redundant\n# \"Compact\" entry:\n#   - paragraphs that resolve to a real style (Casebook Title/Subtitle,\n#     Section Number/Title/Subtitle, Resource Number/Title/Link) keep\n#     that style -- re-applying it collapses the duplicate pStyle down\n#     to a single entry;\n#   - the one paragraph whose *only* style was \"Compact\" is switched to\n#     \"Body Text\" instead.\n\n$d = $word.ActiveDocument\n\n# Styles that, when seen as a paragraph's *current* effective style,\n# mark it as one of the \"Compact + real style\" paragraphs targeted by\n# this edit. Re-assigning the same style collapses the stray duplicate\n# \"Compact\" pStyle entry that precedes it.\n$keepStyles = @(\n    \"Casebook Title\",\n    \"Casebook Subtitle\",\n    \"Section Number\",\n    \"Section Title\",\n    \"Section Subtitle\",\n    \"Resource Number\",\n    \"Resource Title\",\n    \"Resource Link\"\n)\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $styleName = $p.Style.NameLocal\n\n    if ($styleName -eq \"Compact\") {\n        # The lone paragraph whose only pStyle was \"Compact\" becomes \"Body Text\".\n        $p.Style = \"Body Text\"\n    } elseif ($keepStyles -contains $styleName) {\n        $p.Style = $styleName\n    }\n}\n"}
